$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 74 values (quarter 01-01-2021) ---
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 24
$ws.Range("H74").Value = 341

# --- Add new row 75 (quarter 01-04-2021) ---
# Column A holds a text-formatted date label ("01-04-2021"). Assigning the
# string directly would be auto-parsed into a date serial by the General
# number format, so we enter it as a formula returning the literal text and
# then collapse it to a plain value via copy / paste-special-values. This
# keeps the cell on the default (unstyled) format, matching how the other
# date-label cells in the column are stored.
$ws.Range("A75").Formula = "=""01-04-2021"""
$ws.Range("A75").Copy()
$ws.Range("A75").PasteSpecial(-4163)

$ws.Range("B75").Value = 2007
$ws.Range("C75").Value = 15
$ws.Range("D75").Value = -1
$ws.Range("E75").Value = -2
$ws.Range("F75").Value = -17
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 2008
$ws.Range("I75").Value = 21
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 21
